$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8988.88228765323
$ws.Range("F2").Value = -8.13299421102321

$ws.Range("C3").Value = 9772.18316474197
$ws.Range("F3").Value = 263.086472568268

$ws.Range("C4").Value = 9747.01079535469
$ws.Range("F4").Value = 282.416284422915
